$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.182.03"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Value = "2.351.32"
$ws.Range("E3").Value = "  +1.29%  "

$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.72"
$ws.Range("E5").Value = "  +1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.77"
$ws.Range("E6").Value = "  +2.28%  "

$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.56"
$ws.Range("E10").Value = "  +4.74%  "

$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("E12").Value = "  +2.27%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.772.68"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.79"
$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("D15").Value = "58.142.46"
$ws.Range("E15").Value = "  +1.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000133"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").Value = "2.354.81"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.70"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "332.19"
$ws.Range("E19").Value = "  -1.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.27"
$ws.Range("E20").Value = "  +2.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.81"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("E25").Value = "  -2.85%  "

$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.44"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("E30").Value = "  +1.40%  "

$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("E32").Value = "  +12.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.44"
$ws.Range("E33").Value = "  -0.62%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  +6.58%  "

$ws.Range("E36").Value = "  +0.83%  "

$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("E38").Value = "  +4.20%  "

$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "144.95"
$ws.Range("E40").Value = "  -2.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "295.37"
$ws.Range("E41").Value = "  +4.77%  "

$ws.Range("E42").Value = "  +0.80%  "

$ws.Range("E43").Value = "  +0.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0947"
$ws.Range("E44").Value = "  +1.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.19"
$ws.Range("E45").Value = "  +1.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0502"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.54"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("E51").Value = "  +0.44%  "
